$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - match formatting of the existing header cells (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF) - identical values per row
$values = @{
    2  = 7
    3  = 8
    4  = 9
    5  = 4
    6  = 7
    7  = 10
    8  = 9
    9  = 9
    10 = 7
    11 = 8
    12 = 7
    13 = 7
    14 = 7
    15 = 6
    16 = 7
    17 = 6
    18 = 5
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 9).Value = $v
    $ws.Cells.Item($row, 10).Value = $v
}
